# Daily attendance processing - normalize "Recorded By" (column G) ordering.
# For each data row, the comma-separated list of recorder identities in
# column G is reversed (e.g. "System, user@example.com" becomes
# "user@example.com, System"), except for the specific combination
# "admin@admin.com, System", which is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -eq "admin@admin.com, System") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Length -le 1) {
        continue
    }

    $revParts = $parts[($parts.Length - 1)..0]
    $newVal = [string]::Join(", ", $revParts)

    $cell.Value2 = $newVal
}
